$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 95
$ws.Range("I2").Value = 281
$ws.Range("J2").Value = 1220
$ws.Range("K2").Value = 5
$ws.Range("L2").Value = 337
$ws.Range("M2").Value = 16
$ws.Range("N2").Value = 208
$ws.Range("P2").Value = 5
$ws.Range("Q2").Value = 3
$ws.Range("R2").Value = 11
$ws.Range("S2").Value = 135
$ws.Range("T2").Value = 215
$ws.Range("U2").Value = 14
$ws.Range("V2").Value = 1912
$ws.Range("X2").Value = 1859
$ws.Range("Y2").Value = 4
$ws.Range("Z2").Value = 35
$ws.Range("AA2").Value = 15
